# /tmp/work/edit.ps1
# Applies the "Finished module design breakdown (subject to approval)" edit:
#  - Retitle the Command/Control module to just "Control Module"
#  - Rewrite the four module-description paragraphs with the updated copy
#  - Tweak the Normal style font colour from "automatic" to RGB 00000A

$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false,
                                      $true, 1, $false, $newText, 2)
    if (-not $found) {
        throw "Find/Replace did not find expected text: $oldText"
    }
}

# --- "Module Title" heading for the first module ---
$titleOld = "Module Title: Command/Control Module"
$titleNew = "Module Title: Control Module"
Replace-Text $titleOld $titleNew

# --- Body paragraph for the (now) "Control Module" ---
$controlOld = "The Command and Control Module is responsible for the main functions of the device, running the main PID control loop and interfacing with the gyroscope and accelerometer. The program loop is designed to be short enough to allow regular updating of motor speeds through PWM control of the ESCs. The base for this module is the ATMEGA32u4 of the Arduino Leonardo, chosen for its plethora of 4 16-bit timers to allow enough graduation in motor speeds for fine control of the device. This will communicate with the comms module through the UART protocol to receive the user input and will return the motor speed and orientation for logging. The ESCs will take the PWM from the Arduino and produce higher current PWM from the battery to control the 4 motors."
$controlNew = "The Control Module is responsible for the main functions of the device, running the main PID control loop and interfacing with the gyroscope and accelerometer. The program loop is designed to be short enough to allow regular updating of motor speeds through PWM control of the ESCs. The base for this module is the ATMEGA32u4 of the Arduino Leonardo, chosen for the wide availability of sensor interface libraries and enough 16-bit timers to allow enough graduation in motor speeds for fine control of the device. This will interface with the Communications module through the UART protocol to receive the user input and will return the motor speed and device orientation for logging and telemetry. The ESCs isolate the high current power circuitry of the motors from the low current micro-controller, taking in low current PWM with a maximum 10% duty cycle and providing the 4 motors with a PWM signal at much higher currents."
Replace-Text $controlOld $controlNew

# --- Body paragraph for the Sensing Module ---
$sensingOld = "The sensing module of the "
$sensingNew = "The sensing capabilities of the system is based around the MPU6050 gyroscope/accelerometer IMU and a Sharp GP2Y0A41SK0F infra-red proximity sensor. The MPU6050 has an onboard DMP (Digital Motion Processor) which will be used to relieve load on the Control Module Arduino by converting the raw data from the gyroscope into angles for yaw, pitch and roll on chip, before sending this data over an I2C bus ready for use in the stabilisation algorithm. The Infra-red sensor will be mounted on the base of the drone to detect low flight altitudes and semi-automate the landing procedure. The output of this sensor is an analogue voltage which will be fed into an ADC on the Commnunications Module Il Matto board, as this data is not necessary for stable flight. "
Replace-Text $sensingOld $sensingNew

# --- Body paragraph for the Communications Module ---
$commsOld = "The Communications module is based around the ATMEGA644p microcontroller on an ‘Il Matto’ breakout board. It functions as the main communications hub of the system, interfacing with the RFM12B-S2 transceivers over SPI to provide the uplink and downlink to the base station and controller. It will perform some basic processing of the instructions from the controller, passing them through to the command module to create a new setpoint for the controller. This interfacing will be done overt UART, and will receive logging data from the IMU. This logging data will be periodically logged to an SD card through an SPI interface and back to the base station as telemetry information. The communications microcontroller will also be responsible for reading from the IR proximity sensors through its onboard ADCs and controlling the servo-powered cargo hook."
$commsNew = "The Communications module is based around the ATMEGA644p microcontroller on an ‘Il Matto’ breakout board. It functions as the main communications hub of the system, interfacing with the RFM12B-S2 transceivers over SPI to provide the uplink and downlink to the base station and controller. It will perform some basic processing of the instructions from the controller, passing them through to the control module to create a new setpoint for the controller. This interfacing with the command module will be done over UART, and will receive in return logging data from the IMU. This logging data will be periodically written to an SD card sharing the SPI bus and also transmitted back to the base station as telemetry information. The communications microcontroller will also be responsible for reading from the IR proximity sensors through its onboard ADCs and controlling the servo-powered cargo hook through a PWM data signal."
Replace-Text $commsOld $commsNew

# --- Body paragraph for the Base Station / Ground Control module ---
$baseOld = "The Base Station of the system is formed around the core of another Il Matto ATMEGA644p. It will take input from the user through a combination of potentiometer voltages from joysticks being fed into the onboard ADCs and digital inputs from switches and push buttons on the controller. There will also be a UART connection to a host PC to be able to update PID constant values without needing to reprogram the command module, as well as to make the display of telemetry and debug information easier to implement and use. The base station will be connected to a RFM12B-S2 radio transceiver module over SPI bus to allow it to communicate with the quadcopter wirelessly while in flight."
$baseNew = "The Base Station of the system is formed around the core of another Il Matto ATMEGA644p board. It will take input from the user through a combination of joysticks and buttons on the HID (Human Interface Device) controller. The joysticks contain dual potentiometers, which will be directed into the ADCs of the Il Matto to extract values. There will also be a UART connection to a host PC to be able to update PID constant values without needing to reprogram the command module, as well as to make the display of telemetry and debug information easier to implement and use. The base station will be connected to a RFM12B-S2 radio transceiver module over an SPI bus to allow it to communicate with the quadcopter wirelessly while in flight."
Replace-Text $baseOld $baseNew

# --- Normal style: font colour "Automatic" -> RGB(0x00, 0x00, 0x0A) ---
# Word.Font.Color takes a 0xBBGGRR packed long, so RGB 00000A packs to 0x0A0000.
$normalStyle = $d.Styles("Normal")
$normalStyle.Font.Color = 0x0A0000

